$wb = $excel.ActiveWorkbook

# --- sheet2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("C3").Value = 0.9812
$ws.Range("D3").Value = 0.3953
$ws.Range("C5").Value = 0.972
$ws.Range("D5").Value = 0.1296
$ws.Range("C7").Value = 0.9611
$ws.Range("D7").Value = 0.0325
$ws.Range("C9").Value = 0.9539
$ws.Range("D9").Value = 0.0135
$ws.Range("C11").Value = 0.9649
$ws.Range("D11").Value = 0.0526
$ws.Range("C13").Value = 0.9658
$ws.Range("D13").Value = 0.0591
$ws.Range("C15").Value = 0.9673
$ws.Range("D15").Value = 0.0713

# --- sheet3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("C3").Value = 1.0895
$ws.Range("D3").Value = 0.2981
$ws.Range("C4").Value = 0.4149
$ws.Range("D4").Value = 0.5204
$ws.Range("C5").Value = 0.009900000000000001
$ws.Range("D5").Value = 0.9208
$ws.Range("C6").Value = 0.6083
$ws.Range("D6").Value = 0.4365
$ws.Range("C7").Value = 0.5971
$ws.Range("D7").Value = 0.4408
$ws.Range("C8").Value = 0.37
$ws.Range("D8").Value = 0.5438
$ws.Range("C9").Value = 0.0221
$ws.Range("D9").Value = 0.882

# --- sheet4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("E3").Value = 167
$ws.Range("F3").Value = 0.01
$ws.Range("G3").Value = 0.9203
$ws.Range("H3").Value = 0.0001
$ws.Range("E4").Value = 167
$ws.Range("F4").Value = 0.0385
$ws.Range("G4").Value = 0.8447
$ws.Range("H4").Value = 0.0002
$ws.Range("E5").Value = 167
$ws.Range("F5").Value = 0.2799
$ws.Range("G5").Value = 0.5974
$ws.Range("H5").Value = 0.0017
$ws.Range("E6").Value = 167
$ws.Range("F6").Value = 0.1009
$ws.Range("G6").Value = 0.7511
$ws.Range("H6").Value = 0.0005999999999999999
$ws.Range("E7").Value = 167
$ws.Range("F7").Value = 0.06
$ws.Range("G7").Value = 0.8067
$ws.Range("E8").Value = 167
$ws.Range("F8").Value = 0.0755
$ws.Range("G8").Value = 0.7839
$ws.Range("H8").Value = 0.0005
$ws.Range("E9").Value = 167
$ws.Range("F9").Value = 0.08500000000000001
$ws.Range("G9").Value = 0.7711
$ws.Range("H9").Value = 0.0005

# --- sheet5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H3").Value = -0.1027
$ws.Range("I3").Value = 155.358
$ws.Range("K3").Value = 0.9183
$ws.Range("M3").Value = -0.0156
$ws.Range("N3").Value = 0.9183
$ws.Range("H4").Value = 0.1892
$ws.Range("I4").Value = 125.4377
$ws.Range("K4").Value = 0.8502
$ws.Range("M4").Value = 0.0306
$ws.Range("N4").Value = 0.8502
$ws.Range("H5").Value = -0.5328000000000001
$ws.Range("I5").Value = 147.3162
$ws.Range("K5").Value = 0.595
$ws.Range("M5").Value = -0.08260000000000001
$ws.Range("N5").Value = 0.595
$ws.Range("H6").Value = 0.3082
$ws.Range("I6").Value = 128.5122
$ws.Range("K6").Value = 0.7584
$ws.Range("M6").Value = 0.0496
$ws.Range("N6").Value = 0.7584
$ws.Range("H7").Value = -0.2489
$ws.Range("I7").Value = 151.2849
$ws.Range("K7").Value = 0.8038
$ws.Range("M7").Value = -0.0383
$ws.Range("N7").Value = 0.8038
$ws.Range("H8").Value = -0.2689
$ws.Range("I8").Value = 132.9081
$ws.Range("K8").Value = 0.7885
$ws.Range("M8").Value = -0.0429
$ws.Range("N8").Value = 0.7885
$ws.Range("H9").Value = 0.2926
$ws.Range("I9").Value = 145.8583
$ws.Range("K9").Value = 0.7702
$ws.Range("M9").Value = 0.0455
$ws.Range("N9").Value = 0.7702

# L column (p-corr) values are stored as text in this workbook, not numbers.
$ws.Range("L3:L9").NumberFormat = "@"
$ws.Range("L3").Value = "0.17"
$ws.Range("L4").Value = "0.172"
$ws.Range("L5").Value = "0.193"
$ws.Range("L6").Value = "0.177"
$ws.Range("L7").Value = "0.174"
$ws.Range("L8").Value = "0.175"
$ws.Range("L9").Value = "0.176"
